$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''304.55'
$ws.Range('E2').Value = '''1.01%'
$ws.Range('D3').Value = '''35.96'
$ws.Range('E3').Value = '''-3.91%'
$ws.Range('D4').Value = '''5.106'
$ws.Range('E4').Value = '''2.14%'
$ws.Range('D5').Value = '''0.07842'
$ws.Range('E5').Value = '''-0.08%'
$ws.Range('E6').Value = '''-2.55%'
$ws.Range('D7').Value = '''7.938'
$ws.Range('E7').Value = '''-1.16%'
$ws.Range('B8').Value = 'GateToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D8').Value = '''4.106'
$ws.Range('E8').Value = '''2.11%'
$ws.Range('B9').Value = 'MXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D9').Value = '''0.9188'
$ws.Range('E9').Value = '''1.14%'
$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D10').Value = '''0.09695'
$ws.Range('E10').Value = '''1.28%'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').Value = '''0.1855'
$ws.Range('E11').Value = '''-1.90%'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').Value = '''0.08598'
$ws.Range('E12').Value = '''1.38%'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').Value = '''0.03480'
$ws.Range('E13').Value = '''-1.29%'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').Value = '''0.09926'
$ws.Range('E14').Value = '''-0.28%'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').Value = '''0.001443'
$ws.Range('E15').Value = '''-2.71%'
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').Value = '''0.005738'
$ws.Range('E16').Value = '''0.79%'
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').Value = '''3.462'
$ws.Range('E17').Value = '''-0.11%'
$ws.Range('D18').Value = '''2.374'
$ws.Range('E18').Value = '''14.84%'
$ws.Range('D19').Value = '''0.3427'
$ws.Range('E19').Value = '''-1.04%'
$ws.Range('D20').Value = '''0.1302'
$ws.Range('E20').Value = '''-0.50%'
$ws.Range('D21').Value = '''4.819'
$ws.Range('D23').Value = '''0.04527'
$ws.Range('E23').Value = '''-2.73%'
$ws.Range('D24').Value = '''0.005078'
$ws.Range('E24').Value = '''14.08%'
$ws.Range('D25').Value = '''0.001235'
$ws.Range('E25').Value = '''0.53%'
$ws.Range('D26').Value = '''0.0001603'
$ws.Range('E26').Value = '''33.51%'
$ws.Range('D27').Value = '''0.0004758'
$ws.Range('E27').Value = '''0.22%'
$ws.Range('D39').Value = '''0.01832'
$ws.Range('E39').Value = '''4.12%'
$ws.Range('D40').Value = '''0.04716'
$ws.Range('E40').Value = '''-0.60%'
$ws.Range('D41').Value = '''0.007776'
$ws.Range('D42').Value = '''0.1398'
$ws.Range('E42').Value = '''0.42%'
$ws.Range('D43').Value = '''0.007754'
$ws.Range('E43').Value = '''1.09%'
$ws.Range('D44').Value = '''0.002234'
$ws.Range('E44').Value = '''0.13%'
$ws.Range('D45').Value = '''0.01105'
$ws.Range('E45').Value = '''12.50%'
$ws.Range('D46').Value = '''0.00006400'
$ws.Range('E46').Value = '''5.39%'
$ws.Range('D47').Value = '''0.00000000751'
$ws.Range('E47').Value = '''0.24%'
$ws.Range('D48').Value = '''0.0005810'
$ws.Range('E48').Value = '''0.16%'
$ws.Range('D49').Value = '''47.82'
$ws.Range('E49').Value = '''451.54%'
$ws.Range('D50').Value = '''0.002003'
$ws.Range('E50').Value = '''-25.48%'
$ws.Range('D51').Value = '''0.00002104'
$ws.Range('E51').Value = '''0.24%'
